# Generate Report for Handoff
# Updates the localization-status workbook from a "handed back" snapshot to a
# freshly generated "ready for handoff" snapshot: new source/target file
# identifiers, refreshed handoff timestamp, cleared (not-yet-populated)
# "Latest Target File" / "Latest Handback File" columns, and a reset/blank
# "Latest Handback DateTime".

$wb = $excel.ActiveWorkbook

$oldMd1 = "50527144-a073-47db-9c7c-0e38a0676b0f.md"
$oldMd2 = "a45f4c37-42f6-490b-8d01-a84c223ce2ca.md"
$newMd1 = "437cfbd5-767f-4178-a01b-f91116985aef.md"
$newMd2 = "ffff59a32bda-caa8-4d76-8caa-8e980c81ed1f.md"

$newStatus = "Ready for handoff"
$newHandoffDate = "2016-03-24 10:18:19"
$newHandoffDatetime = "2016-03-24 10:18:15"
$newHandbackDatetime = "0001-01-01 00:00:00"

$newXlfZh = "437cfbd5-767f-4178-a01b-f91116985aef.d8c2c853bf6105a22b55831bbb4018e9b23b0fbb.zh-cn.xlf"
$newXlfDe = "437cfbd5-767f-4178-a01b-f91116985aef.d8c2c853bf6105a22b55831bbb4018e9b23b0fbb.de-de.xlf"

$srcRepoBase = "https://github.com/OpenLocalizationTest/oltest/blob/7c7e99c95fcfbb66d19cc97419ef4d5077845f78/e2e"
$zhXlfBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d1b2332d28c71cab1063bad17201de843e8e685a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht"
$deXlfBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ccf8334df655560a5a02062db52de20bac5ab217/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht"

$missing = [System.Reflection.Missing]::Value

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newMd1
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("D2").Value = $newHandoffDate

$wsOverview.Range("A3").Value = $newMd2
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus
$wsOverview.Range("D3").Value = $newHandoffDate

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "$srcRepoBase/$newMd1", $missing, $missing, $newMd1)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "$srcRepoBase/$newMd2", $missing, $missing, $newMd2)

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newMd1
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("D2").Value = $newXlfZh
$wsZh.Range("E2").Value = $newHandoffDatetime
$wsZh.Range("H2").Value = $newHandbackDatetime
$wsZh.Range("F2:G2").Clear()

$wsZh.Range("A3").Value = $newMd2
$wsZh.Range("C3").Value = $newStatus
$wsZh.Range("D3").Value = $newXlfZh
$wsZh.Range("E3").Value = $newHandoffDatetime
$wsZh.Range("H3").Value = $newHandbackDatetime
$wsZh.Range("F3:G3").Clear()

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "$srcRepoBase/$newMd1", $missing, $missing, $newMd1)
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), "$zhXlfBase/$newXlfZh", $missing, $missing, $newXlfZh)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "$srcRepoBase/$newMd2", $missing, $missing, $newMd2)
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), "$zhXlfBase/$newXlfZh", $missing, $missing, $newXlfZh)

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newMd1
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("D2").Value = $newXlfDe
$wsDe.Range("E2").Value = $newHandoffDatetime
$wsDe.Range("H2").Value = $newHandbackDatetime
$wsDe.Range("F2:G2").Clear()

$wsDe.Range("A3").Value = $newMd2
$wsDe.Range("C3").Value = $newStatus
$wsDe.Range("D3").Value = $newXlfDe
$wsDe.Range("E3").Value = $newHandoffDatetime
$wsDe.Range("H3").Value = $newHandbackDatetime
$wsDe.Range("F3:G3").Clear()

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "$srcRepoBase/$newMd1", $missing, $missing, $newMd1)
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), "$deXlfBase/$newXlfDe", $missing, $missing, $newXlfDe)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "$srcRepoBase/$newMd2", $missing, $missing, $newMd2)
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), "$deXlfBase/$newXlfDe", $missing, $missing, $newXlfDe)
